# Automatische test-sync: 2025-08-13 22:51:50
# Appends a new log row (row 30) to the "Logs" sheet, mirroring the most
# recent "Demo inplannen" entry, and bumps the matching tally on the
# "Dashboard" sheet from 28 to 29.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 30

$logs.Cells.Item($newRow, 1).Value = "Demo inplannen"
$logs.Cells.Item($newRow, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-13 22:51:45"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Update the category tally on the Dashboard sheet.
$dashboard.Cells.Item(2, 2).Value = 29

# Extend the conditional-formatting ranges to cover the new row.
$logs.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D30"))
$logs.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G30"))
$logs.Range("H2:H29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H30"))
$logs.Range("I2:I29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I30"))
$logs.Range("J2:J29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J30"))
